# Update Italy commercial mapping and area to match CR typologies in 2011 census
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: apply the CR-typology relabeling to the first 8 "CR/LFINF" lines of a
# mapping scheme block: the first 4 (CDL damage-level rows) become CDL3, and the
# next 4 (CDN damage-level rows) become CDL2. Remaining lines are unchanged.
function Update-CrTypology($text) {
    $lines = $text.Split("`n")
    for ($i = 0; $i -lt $lines.Length; $i++) {
        if ($i -lt 4) {
            $lines[$i] = $lines[$i].Replace("CDL/", "CDL3/")
        } elseif ($i -lt 8) {
            $lines[$i] = $lines[$i].Replace("CDN/", "CDL2/")
        }
    }
    return [string]::Join("`n", $lines)
}

# Read the existing mapping scheme text for each occupancy column.
$officesText = $ws.Range("B2").Value()
$tradeText   = $ws.Range("C2").Value()
$hotelsText  = $ws.Range("D2").Value()

# Apply the 2011-census CR typology relabeling to each block.
$officesText = Update-CrTypology $officesText
$tradeText   = Update-CrTypology $tradeText
$hotelsText  = Update-CrTypology $hotelsText

# Write the relabeled text back into its original column (Offices stays in B2,
# Trade stays in C2, Hotels stays in D2).
$ws.Range("D2").Value = $hotelsText
$ws.Range("C2").Value = $tradeText
$ws.Range("B2").Value = $officesText

# Update the active selection to match the saved view state.
[void]$ws.Range("B3").Select()
